# Ephyz unit tests in DB
# Adds three new rows (21-23) to Sheet1 documenting GC ADP / half-duration
# measurements from Stroh et. al. (2012), mirroring the existing table
# layout/style used by the other paper rows (e.g. row 18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: paper-level summary row for the GC ADP/LLD entry ---
$ws.Range("A21").Value = "Stroh et. al. (2012)"
$ws.Range("B21").Value = "NMDA Receptor-Dependent Synaptic Activation of TRPC`nChannels in Olfactory Bulb Granule Cells"
$ws.Range("B21").WrapText = $false
$ws.Range("C21").Value = "Mice"
$ws.Range("D21").Value = "GC"
$ws.Range("E21").Value = "GCs show ADPs in vivo to somatic stimulation, and LLDs to glomerular stimulation"

# --- Row 22: ADP amplitude measurement row ---
$ws.Range("A22").Value = "Stroh et. al. (2012)"
$ws.Range("B22").Value = "NMDA Receptor-Dependent Synaptic Activation of TRPC`nChannels in Olfactory Bulb Granule Cells"
$ws.Range("B22").WrapText = $false
$ws.Range("C22").Value = "Mice"
$ws.Range("D22").Value = "GC"
$ws.Range("E22").Value = "ADP amplitude"
$ws.Range("F22").Value = "11.1+-4.7 mV"
$ws.Range("G22").Value = 49
$ws.Range("H22").Value = 21
$ws.Range("I22").Value = "not REPORTED"
$ws.Range("J22").Value = "Soma injections 1000pa for 1ms. a mean sAP-ADP amplitude above resting potential"
$ws.Range("K22").Value = "STD"

# --- Row 23: ADP half-duration measurement row ---
$ws.Range("E23").Value = "ADP half-duration"
$ws.Range("F23").Value = "42+-22 ms"
$ws.Range("G23").Value = 49
$ws.Range("H23").Value = 21
$ws.Range("I23").Value = "not REPORTED"
$ws.Range("J23").Value = "Soma injections 1000pa for 1ms"
$ws.Range("K23").Value = "Throughout the paper, tau1/2 denotes halfdurations from the peak amplitude of Vm onward, measured between the onset of the afterdepolarization (ADP) right after the sodium spike and one-half of its maximum amplitude"
$ws.Range("K23").WrapText = $false

# Match the existing table's style on the paper-title cells (same look as
# B18, which already carries this "applyAlignment" style in the template).
$ws.Range("B18").ShrinkToFit = $false
$ws.Range("B21").ShrinkToFit = $false
$ws.Range("B22").ShrinkToFit = $false
$ws.Range("K23").ShrinkToFit = $false

# Move the active selection down, the same way Excel leaves it positioned
# one row below the last-entered row after data entry.
$ws.Range("E20").Select()
